$wb = $excel.ActiveWorkbook

# --- Worksheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 846.73334
$ws.Range("I19").Value = 589.4
$ws.Range("K19").Value = 589.4
$ws.Range("M19").Value = -414.4
$ws.Range("H33").Value = 5263406.5
$ws.Range("I33").Value = 8333641
$ws.Range("K33").Value = 8333641
$ws.Range("M33").Value = -8333412
$ws.Range("H40").Value = 2457.8572
$ws.Range("I40").Value = 1800.75
$ws.Range("J40").Value = 3334
$ws.Range("K40").Value = 1800.75
$ws.Range("L40").Value = 3334
$ws.Range("M40").Value = -1625.75
$ws.Range("N40").Value = -3684
$ws.Range("H88").Value = 2454.2727
$ws.Range("I88").Value = 5251.5
$ws.Range("J88").Value = 1832.6666
$ws.Range("K88").Value = 5251.5
$ws.Range("L88").Value = 1832.6666
$ws.Range("M88").Value = -4845.5
$ws.Range("N88").Value = -2644.6666
$ws.Range("H91").Value = 2454.2727
$ws.Range("I91").Value = 5251.5
$ws.Range("J91").Value = 1832.6666
$ws.Range("K91").Value = 5251.5
$ws.Range("L91").Value = 1832.6666
$ws.Range("M91").Value = -3847.5
$ws.Range("N91").Value = -4640.6666
$ws.Range("H137").Value = 5560913
$ws.Range("I137").Value = 12503015
$ws.Range("K137").Value = 37509045
$ws.Range("M137").Value = -37506495
$ws.Range("H138").Value = 5013.364
$ws.Range("I138").Value = 2925
$ws.Range("J138").Value = 5477.4443
$ws.Range("K138").Value = 8775
$ws.Range("L138").Value = 16432.3329
$ws.Range("M138").Value = -3635
$ws.Range("N138").Value = -26712.3329
$ws.Range("H141").Value = 738.2
$ws.Range("I141").Value = 647.5
$ws.Range("K141").Value = 1942.5
$ws.Range("M141").Value = 3237.5

# --- Worksheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3316.9807
$ws.Range("I32").Value = 3084.279
$ws.Range("K32").Value = 3084.279
$ws.Range("M32").Value = -2797.279
$ws.Range("H44").Value = 93372.5
$ws.Range("J44").Value = 93372.5
$ws.Range("L44").Value = 93372.5
$ws.Range("N44").Value = -94348.5
$ws.Range("H45").Value = 1104.4166
$ws.Range("I45").Value = 1095.9
$ws.Range("J45").Value = 1147
$ws.Range("K45").Value = 1095.9
$ws.Range("L45").Value = 1147
$ws.Range("M45").Value = -718.9000000000001
$ws.Range("N45").Value = -1901
$ws.Range("H55").Value = 32499.25
$ws.Range("J55").Value = 59900
$ws.Range("L55").Value = 59900
$ws.Range("N55").Value = -60530
$ws.Range("H69").Value = 503448.28
$ws.Range("J69").Value = 503448.28
$ws.Range("L69").Value = 503448.28
$ws.Range("N69").Value = -504946.28
$ws.Range("H72").Value = 503448.28
$ws.Range("J72").Value = 503448.28
$ws.Range("L72").Value = 1510344.84
$ws.Range("N72").Value = -1517832.84

# --- Worksheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1650
$ws.Range("I99").Value = 800
$ws.Range("J99").Value = 2500
$ws.Range("K99").Value = 800
$ws.Range("L99").Value = 2500
$ws.Range("M99").Value = 698
$ws.Range("N99").Value = -5496

# --- Worksheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2087395
$ws.Range("I31").Value = 1195.0625
$ws.Range("J31").Value = 3130495
$ws.Range("K31").Value = 1195.0625
$ws.Range("L31").Value = 3130495
$ws.Range("M31").Value = -900.0625
$ws.Range("N31").Value = -3131085
$ws.Range("H34").Value = 2087395
$ws.Range("I34").Value = 1195.0625
$ws.Range("J34").Value = 3130495
$ws.Range("K34").Value = 1195.0625
$ws.Range("L34").Value = 3130495
$ws.Range("M34").Value = -993.0625
$ws.Range("N34").Value = -3130899
$ws.Range("H53").Value = 54999.5
$ws.Range("J53").Value = 54999.5
$ws.Range("L53").Value = 54999.5
$ws.Range("N53").Value = -56213.5
$ws.Range("H62").Value = 24999.625
$ws.Range("I62").Value = 19998
$ws.Range("J62").Value = 25714.143
$ws.Range("K62").Value = 19998
$ws.Range("L62").Value = 25714.143
$ws.Range("M62").Value = -19374
$ws.Range("N62").Value = -26962.143
$ws.Range("H65").Value = 24999.625
$ws.Range("I65").Value = 19998
$ws.Range("J65").Value = 25714.143
$ws.Range("K65").Value = 99990
$ws.Range("L65").Value = 128570.715
$ws.Range("M65").Value = -96870
$ws.Range("N65").Value = -134810.715
$ws.Range("H87").Value = 45000
$ws.Range("I87").Value = 45000
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 45000
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -43814
$ws.Range("N87").ClearContents()
$ws.Range("H90").Value = 45000
$ws.Range("I90").Value = 45000
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 135000
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -129072
$ws.Range("N90").ClearContents()
$ws.Range("H94").Value = 4426
$ws.Range("J94").Value = 4426
$ws.Range("L94").Value = 4426
$ws.Range("N94").Value = -5328
$ws.Range("H122").Value = 999.75
$ws.Range("I122").Value = 999.6667
$ws.Range("K122").Value = 2999.0001
$ws.Range("M122").Value = -549.0001000000002

# --- Worksheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H52").Value = 3000
$ws.Range("I52").Value = 1000
$ws.Range("J52").Value = 5000
$ws.Range("K52").Value = 3000
$ws.Range("L52").Value = 15000
$ws.Range("M52").Value = -2734
$ws.Range("N52").Value = -15532
$ws.Range("H68").Value = 120538.12
$ws.Range("I68").Value = 334190.84
$ws.Range("K68").Value = 1002572.52
$ws.Range("M68").Value = -1001761.52
$ws.Range("H71").Value = 120538.12
$ws.Range("I71").Value = 334190.84
$ws.Range("K71").Value = 3007717.56
$ws.Range("M71").Value = -3003661.56
$ws.Range("H94").Value = 3824.6667
$ws.Range("I94").Value = 3482.6667
$ws.Range("J94").Value = 4166.6665
$ws.Range("K94").Value = 10448.0001
$ws.Range("L94").Value = 12499.9995
$ws.Range("M94").Value = -9772.000100000001
$ws.Range("N94").Value = -13851.9995
$ws.Range("H137").Value = 1985.6842
$ws.Range("I137").Value = 1208
$ws.Range("J137").Value = 2551.2727
$ws.Range("K137").Value = 3624
$ws.Range("L137").Value = 7653.8181
$ws.Range("M137").Value = 1476
$ws.Range("N137").Value = -17853.8181
$ws.Range("H138").Value = 696339.3
$ws.Range("I138").Value = 760.5714
$ws.Range("J138").Value = 2319356.2
$ws.Range("K138").Value = 2281.7142
$ws.Range("L138").Value = 6958068.600000001
$ws.Range("M138").Value = 2858.2858
$ws.Range("N138").Value = -6968348.600000001

# --- Worksheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H57").Value = 49721.25
$ws.Range("J57").Value = 49721.25
$ws.Range("L57").Value = 49721.25
$ws.Range("N57").Value = -51361.25
$ws.Range("H123").Value = 39053.43
$ws.Range("J123").Value = 39053.43
$ws.Range("L123").Value = 39053.43
$ws.Range("N123").Value = -43953.43

# --- Worksheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 17993.47
$ws.Range("I7").Value = 18407.416
$ws.Range("K7").Value = 18407.416
$ws.Range("M7").Value = -18295.416
$ws.Range("H46").Value = 5257.864
$ws.Range("I46").Value = 1747.2
$ws.Range("J46").Value = 6290.4116
$ws.Range("K46").Value = 1747.2
$ws.Range("L46").Value = 6290.4116
$ws.Range("M46").Value = -1559.2
$ws.Range("N46").Value = -6666.4116
$ws.Range("H55").Value = 663.3
$ws.Range("J55").Value = 1090.75
$ws.Range("L55").Value = 1090.75
$ws.Range("N55").Value = -1436.75
$ws.Range("H126").Value = 17993.47
$ws.Range("I126").Value = 18407.416
$ws.Range("K126").Value = 55222.24800000001
$ws.Range("M126").Value = -52752.24800000001
$ws.Range("H132").Value = 4071.5476
$ws.Range("I132").Value = 3556.05
$ws.Range("K132").Value = 10668.15
$ws.Range("M132").Value = -8138.150000000001

# --- Worksheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 2527.7273
$ws.Range("I100").Value = 3371.375
$ws.Range("K100").Value = 6742.75
$ws.Range("M100").Value = -6201.75
$ws.Range("H113").Value = 928.9524
$ws.Range("I113").Value = 562.5333000000001
$ws.Range("J113").Value = 1845
$ws.Range("K113").Value = 1687.5999
$ws.Range("L113").Value = 5535
$ws.Range("M113").Value = 482.4000999999998
$ws.Range("N113").Value = -9875
$ws.Range("H132").Value = 5427.4565
$ws.Range("I132").Value = 3773.8696
$ws.Range("J132").Value = 7081.0435
$ws.Range("K132").Value = 11321.6088
$ws.Range("L132").Value = 21243.1305
$ws.Range("M132").Value = -8791.6088
$ws.Range("N132").Value = -26303.1305
$ws.Range("H136").Value = 4844.933
$ws.Range("I136").Value = 3953.2942
$ws.Range("J136").Value = 6010.923
$ws.Range("K136").Value = 11859.8826
$ws.Range("L136").Value = 18032.769
$ws.Range("M136").Value = -9309.882599999999
$ws.Range("N136").Value = -23132.769
